# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly scraped totals (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2532
$ws1.Range("F5").Value = 1433
$ws1.Range("F6").Value = 1114
$ws1.Range("F13").Value = 8809
$ws1.Range("F14").Value = 381
$ws1.Range("F15").Value = 2494
$ws1.Range("F25").Value = 2115
$ws1.Range("F27").Value = 1804
$ws1.Range("F31").Value = 223
$ws1.Range("F33").Value = 116
$ws1.Range("F35").Value = 14
$ws1.Range("F36").Value = 313
$ws1.Range("F38").Value = 263
$ws1.Range("F39").Value = 445
$ws1.Range("F40").Value = 767
$ws1.Range("F42").Value = 272

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2532
$ws4.Range("F5").Value = 1433
$ws4.Range("F7").Value = 1114
$ws4.Range("F14").Value = 8809
$ws4.Range("F15").Value = 381
$ws4.Range("F16").Value = 2494
$ws4.Range("F27").Value = 2115
$ws4.Range("F29").Value = 1804
$ws4.Range("F33").Value = 223
$ws4.Range("F35").Value = 116
$ws4.Range("F37").Value = 14
$ws4.Range("F38").Value = 313
$ws4.Range("F40").Value = 263
$ws4.Range("F41").Value = 445
$ws4.Range("F46").Value = 767
$ws4.Range("F49").Value = 272
